$wb = $excel.ActiveWorkbook

# 1. Fix typo "excedes" -> "exceeds" in the Criterion % Definitions sheet
#    (do this before the rename, while we can still look the sheet up by its old name)
$critWs = $wb.Worksheets.Item("Criterion % Definitions")
$critWs.Cells.Item(2,2).Value = "The percentage of occupied hours where delta T equals or exceeds the threshold (1 kelvin) over the total occupied hours."

# 2. Rename the sheet "Criterion % Definitions" -> "Criterion Definitions"
$critWs.Name = "Criterion Definitions"

# 3. Update the recorded analysis timestamp on the Project Information sheet
$projWs = $wb.Worksheets.Item("Project Information")
$projWs.Cells.Item(11,2).Value = "2022-03-02 16:50:29.101876"

# 4. On the readme sheet, reorder the metadata table columns from
#    index, JobNo, Date, Author, sheet_name
#    to
#    index, JobNo, sheet_name, Date, Author
$ws = $wb.Worksheets.Item("readme")

# Header row - plain text values, safe to set directly (keeps the ListObject/table in sync)
$ws.Cells.Item(1,3).Value = "sheet_name"
$ws.Cells.Item(1,4).Value = "Date"
$ws.Cells.Item(1,5).Value = "Author"

# Data rows 2-12 - use Copy (via helper columns) instead of Value assignment so that
# text that looks numeric (e.g. "20220302") keeps its text type/style instead of being
# coerced into a number.
for ($r=2; $r -le 12; $r++) {
    $ws.Cells.Item($r,5).Copy($ws.Cells.Item($r,7))   # G = old E (sheet_name)
    $ws.Cells.Item($r,3).Copy($ws.Cells.Item($r,8))   # H = old C (Date)
    $ws.Cells.Item($r,4).Copy($ws.Cells.Item($r,9))   # I = old D (Author)
}
for ($r=2; $r -le 12; $r++) {
    $ws.Cells.Item($r,7).Copy($ws.Cells.Item($r,3))   # C (sheet_name) = old E
    $ws.Cells.Item($r,8).Copy($ws.Cells.Item($r,4))   # D (Date) = old C
    $ws.Cells.Item($r,9).Copy($ws.Cells.Item($r,5))   # E (Author) = old D
}
$ws.Range("G2:I12").Clear()

# The sheet_name column for row 3 (index 1) refers to the sheet we just renamed;
# update its cached text to match the new sheet name.
$ws.Cells.Item(3,3).Value = "Criterion Definitions"
